$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '42.509.63'
$ws.Range('E2').Value = '  -2.56%  '
$ws.Range('D3').Value = '2.224.28'
$ws.Range('E3').Value = '  -2.20%  '
$ws.Range('E4').Value = '  +0.34%  '
$ws.Range('D5').Value = '''112.26'
$ws.Range('E5').Value = '  -8.08%  '
$ws.Range('D6').Value = '''297.75'
$ws.Range('E6').Value = '  +12.31%  '
$ws.Range('E7').Value = '  -1.54%  '
$ws.Range('E8').Value = '  +0.00%  '
$ws.Range('D9').Value = '''0.612'
$ws.Range('E9').Value = '  -1.68%  '
$ws.Range('D10').Value = '''45.34'
$ws.Range('E10').Value = '  -5.45%  '
$ws.Range('E11').Value = '  -1.51%  '
$ws.Range('D12').Value = '''54.77'
$ws.Range('E12').Value = '  +0.66%  '
$ws.Range('D13').Value = '''8.90'
$ws.Range('E13').Value = '  -0.81%  '
$ws.Range('E14').Value = '  -2.88%  '
$ws.Range('D15').Value = '''0.952'
$ws.Range('E15').Value = '  +6.86%  '
$ws.Range('D16').Value = '''15.17'
$ws.Range('E16').Value = '  -1.37%  '
$ws.Range('D17').Value = '2.560.08'
$ws.Range('E17').Value = '  -2.18%  '
$ws.Range('D18').Value = '2.262.31'
$ws.Range('E18').Value = '  -0.52%  '
$ws.Range('D19').Value = '42.410.32'
$ws.Range('E19').Value = '  -2.71%  '
$ws.Range('D20').Value = '''7.39'
$ws.Range('E20').Value = '  +6.26%  '
$ws.Range('E21').Value = '  -2.59%  '
$ws.Range('D22').Value = '''73.84'
$ws.Range('E22').Value = '  +2.35%  '
$ws.Range('D23').Value = '''3.54'
$ws.Range('E23').Value = '  +24.26%  '
$ws.Range('D24').Value = '''2.30'
$ws.Range('E24').Value = '  -5.23%  '
$ws.Range('D25').Value = '''229.88'
$ws.Range('E25').Value = '  -2.13%  '
$ws.Range('D26').Value = '''9.46'
$ws.Range('E26').Value = '  -0.41%  '
$ws.Range('D27').Value = '''11.77'
$ws.Range('E27').Value = '  +0.35%  '
$ws.Range('E28').Value = '  -1.50%  '
$ws.Range('D29').Value = '''3.90'
$ws.Range('E29').Value = '  -1.51%  '
$ws.Range('D30').Value = '''38.98'
$ws.Range('E30').Value = '  -6.68%  '
$ws.Range('E31').Value = '  -1.57%  '
$ws.Range('E32').Value = '  -3.52%  '
$ws.Range('D33').Value = '''174.00'
$ws.Range('E33').Value = '  +1.36%  '
$ws.Range('D34').Value = '''21.17'
$ws.Range('E34').Value = '  -2.06%  '
$ws.Range('D35').Value = '''0.0890'
$ws.Range('E35').Value = '  -1.67%  '
$ws.Range('D36').Value = '''5.77'
$ws.Range('E36').Value = '  +0.82%  '
$ws.Range('B37').Value = 'RenderToken'
$ws.Range('C37').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D37').Value = '''4.90'
$ws.Range('E37').Value = '  +5.32%  '
$ws.Range('B38').Value = 'NEARProtocol'
$ws.Range('C38').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D38').Value = '''4.29'
$ws.Range('E38').Value = '  +5.78%  '
$ws.Range('E39').Value = '  -1.71%  '
$ws.Range('D40').Value = '''0.0369'
$ws.Range('E40').Value = '  -2.40%  '
$ws.Range('E41').Value = '  -3.14%  '
$ws.Range('D42').Value = '''2.50'
$ws.Range('E42').Value = '  -1.16%  '
$ws.Range('D43').Value = '''0.239'
$ws.Range('E43').Value = '  +0.68%  '
$ws.Range('D44').Value = '''71.23'
$ws.Range('E44').Value = '  -5.17%  '
$ws.Range('D45').Value = '''13.17'
$ws.Range('E45').Value = '  -5.50%  '
$ws.Range('E46').Value = '  +0.30%  '
$ws.Range('E47').Value = '  -2.22%  '
$ws.Range('E48').Value = '  -3.93%  '
$ws.Range('E49').Value = '  +4.39%  '
$ws.Range('D50').Value = '''105.20'
$ws.Range('E50').Value = '  +3.98%  '
$ws.Range('D51').Value = '''8.57'
$ws.Range('E51').Value = '  +0.44%  '
